$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")

# Update the "Updated:" banner date
$ws.Range("B2").Value = "Updated: 2023 - 03 - 07"

# Project #9 (Pre-process / online community relation) is now Closed
$ws.Range("E12").Value = "Closed"

# New project #10 - openCV image stitching review
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = "Stitching camera 2D images"
$ws.Range("D13").Value = "No meeting"
$ws.Range("E13").Value = "Open"
$ws.Range("F13").Value = "-"
$ws.Range("G13").Value = "2023 - 03 - 07"

$ws.Range("G7").Select()
